$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$rowUpdates = @{
    2  = 134
    3  = 49
    7  = 1245
    8  = 1528
    9  = 337
    10 = 387
    12 = 145
    15 = 105
    17 = 298
    18 = 321
    19 = 1722
    23 = 662
    25 = 331
    26 = 4144
    27 = 12
    28 = 263
    30 = 481
    32 = 512
    34 = 232
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $rowUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowUpdates[$row]
    }
}
